$p = $ppt.ActivePresentation

# Keep only the first 3 slides (MINT intro, agenda, "Introduction to MINT");
# remove the remaining slides (the history/applications/links deck tail),
# deleting from the end so indices of the slides we keep stay stable.
for ($i = $p.Slides.Count; $i -ge 4; $i--) {
    $p.Slides.Item($i).Delete()
}
